# Trading update: 2026-02-17 08:39:05
# Append a new (OPEN) trade row (row 47 / Trade # 46) to both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

$newRow = 47

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$newRow").Value = 46

    # Date/Time are plain text in this sheet (e.g. "2026-02-17"), not real
    # dates - a leading "'" keeps Excel from auto-converting them to date
    # serials, and resetting the style afterwards drops the quote-prefix
    # formatting flag that COM would otherwise stamp on the cell.
    $ws.Range("B$newRow").Value = "'2026-02-17"
    $ws.Range("B$newRow").Style = "Normal"

    $ws.Range("C$newRow").Value = "'08:39:00"
    $ws.Range("C$newRow").Style = "Normal"

    $ws.Range("D$newRow").Value = "MarketMaking"
    $ws.Range("E$newRow").Value = "UP"
    $ws.Range("F$newRow").Value = 0.86

    # Exit Price is blank (trade still OPEN) but present as an empty text
    # cell, same as the rest of the row's empty cells - "'" produces an
    # empty-string text value once the quote prefix is stripped.
    $ws.Range("G$newRow").Value = "'"
    $ws.Range("G$newRow").Style = "Normal"

    $ws.Range("H$newRow").Value = "OPEN"
    $ws.Range("I$newRow").Value = 0
    $ws.Range("J$newRow").Value = 0
    $ws.Range("K$newRow").Value = 99.53598934440596
    $ws.Range("L$newRow").Value = 0
    $ws.Range("M$newRow").Value = 0
    $ws.Range("N$newRow").Value = 0.6
    $ws.Range("O$newRow").Value = "Normal spread capture: 19600 bps"

    $ws.Range("P$newRow").Value = "'"
    $ws.Range("P$newRow").Style = "Normal"

    $ws.Range("Q$newRow").Value = 0
}
